$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "36.466.20"
$ws.Range("E2").Value = "  -1.49%  "
$ws.Range("D3").Value = "2.036.71"
$ws.Range("E3").Value = "  +2.29%  "
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "230.39"
$ws.Range("E5").Value = "  -12.60%  "
$ws.Range("E6").Value = "  -1.36%  "
$ws.Range("E7").Value = "  -0.12%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "54.79"
$ws.Range("E8").Value = "  -1.25%  "
$ws.Range("E9").Value = "  -1.25%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "56.86"
$ws.Range("E10").Value = "  +1.27%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0748"
$ws.Range("E11").Value = "  -1.54%  "
$ws.Range("E12").Value = "  +0.12%  "
$ws.Range("D13").Value = "2.337.74"
$ws.Range("E13").Value = "  +2.44%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "14.27"
$ws.Range("E14").Value = "  +1.22%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "19.94"
$ws.Range("E15").Value = "  -9.00%  "
$ws.Range("E16").Value = "  -1.43%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.09"
$ws.Range("E17").Value = "  -1.18%  "
$ws.Range("D18").Value = "2.036.47"
$ws.Range("E18").Value = "  +2.50%  "
$ws.Range("D19").Value = "36.583.22"
$ws.Range("E19").Value = "  -0.72%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "5.82"
$ws.Range("E20").Value = "  +15.41%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "67.27"
$ws.Range("E21").Value = "  -3.22%  "
$ws.Range("E22").Value = "  -3.56%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "219.82"
$ws.Range("E23").Value = "  -5.90%  "
$ws.Range("E24").Value = "  +0.01%  "
$ws.Range("E25").Value = "  +1.65%  "
$ws.Range("E26").Value = "  -8.67%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "162.63"
$ws.Range("E27").Value = "  -1.65%  "
$ws.Range("E28").Value = "  -1.36%  "
$ws.Range("E29").Value = "  -0.13%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "18.82"
$ws.Range("E30").Value = "  -2.27%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.35"
$ws.Range("E31").Value = "  +3.91%  "
$ws.Range("E32").Value = "  -1.73%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.34"
$ws.Range("E33").Value = "  -3.81%  "
$ws.Range("E34").Value = "  -2.92%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.47"
$ws.Range("E35").Value = "  +3.20%  "
$ws.Range("E36").Value = "  -2.03%  "
$ws.Range("E37").Value = "  -0.12%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.75"
$ws.Range("E38").Value = "  -2.61%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.76"
$ws.Range("E39").Value = "  +8.15%  "
$ws.Range("E40").Value = "  -6.93%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.53"
$ws.Range("E41").Value = "  +46.94%  "
$ws.Range("E42").Value = "  -3.94%  "
$ws.Range("D43").Value = "1.481.14"
$ws.Range("E43").Value = "  +3.41%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0934"
$ws.Range("E44").Value = "  +2.79%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "92.96"
$ws.Range("E45").Value = "  +4.68%  "
$ws.Range("E46").Value = "  -1.45%  "
$ws.Range("E47").Value = "  -4.74%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "15.49"
$ws.Range("E48").Value = "  +0.72%  "
$ws.Range("E49").Value = "  -1.72%  "
$ws.Range("E50").Value = "  -0.05%  "
$ws.Range("E51").Value = "  +2.18%  "
